# "edit logo & add filtro tipo de fase"
# Fills the newly-introduced "Video" column (I) for every data row with the
# match video link, widens column I to fit the URL, and leaves the sheet's
# selection parked on the filled range (I3:I96) as the author's session did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$videoUrl = "https://youtu.be/7N6VSnEXyFQ"

# Data rows run from 2 to 96 (row 1 is the header, which already has the
# "Video" header label in I1). Every row gets the same shared video link.
$ws.Range("I2:I96").Value = $videoUrl

# Widen column I (bestFit-style) now that it holds the long URL text.
$ws.Columns.Item(9).ColumnWidth = 26.877604166666668

# Match the author's final on-screen selection: I3 active, I3:I96 selected.
$ws.Range("I3:I96").Select() | Out-Null
